$wb = $excel.ActiveWorkbook

$rows = @(
    @{ Sheet = 1; Time = "2025-03-07 20:42:06"; C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x14,0x41,0x0c,"; E = "0x d";  G = "568631262647113770877196"; I = 13 },
    @{ Sheet = 2; Time = "2025-03-07 20:29:35"; C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x15,0x41,0x0c,"; E = "0x e";  G = "568631262647113770942732"; I = 14 },
    @{ Sheet = 3; Time = "2025-03-07 20:51:45"; C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x06,0x41,0x0c,"; E = "0xff"; G = "568631262647113769959692"; I = 255 },
    @{ Sheet = 4; Time = "2025-03-07 20:41:15"; C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x0b,0x40,0x0c,"; E = "0x 3"; G = "568631262647113769959692"; I = 3 }
)

foreach ($row in $rows) {
    $ws = $wb.Worksheets.Item($row.Sheet)

    $ws.Cells.Item(85, 1).Value = $row.Time
    $ws.Cells.Item(85, 2).Value = "0x01,0x90 "
    $ws.Cells.Item(85, 3).Value = $row.C
    $ws.Cells.Item(85, 4).Value = "0x01,0x90,"
    $ws.Cells.Item(85, 5).Value = $row.E
    $ws.Cells.Item(85, 6).Value = 400
    $ws.Cells.Item(85, 7).NumberFormat = "@"
    $ws.Cells.Item(85, 7).Value = $row.G
    $ws.Cells.Item(85, 8).Value = 400
    $ws.Cells.Item(85, 9).Value = $row.I
}
